$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.3632265282951721,
    -0.087979172526694072,
    -0.0089999997813485777,
    -0.011999999938343109,
    -0.005999999786560295,
    -0.0059999997816078121,
    -0.019999999735720309,
    -0.019999999735101248,
    -0.005999999780666343,
    -0.0059999997810038508,
    -0.0044999997860095675,
    -0.0059999997811823746,
    0.062331131190581424,
    -0.01199999976241628,
    -0.0059999997811521766,
    -0.0059999997804043304,
    -0.0059999997794530913,
    -0.0089999997694336642,
    -0.008999999786910351,
    -0.0089999997849705693,
    -0.0089999997846543778,
    -0.0089999997844190105,
    -0.0089999997792862274,
    -0.041999999668226096,
    -0.041999999666370691,
    -0.0059999997810358252,
    -0.0059999997801130078,
    0.019781416904098847,
    -0.011999999754067403,
    -0.019999999726365125,
    -0.01499999974028654,
    -0.020999999720375584,
    -0.0059999997690614038
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
